# Changed local db as docker db
$wb = $excel.ActiveWorkbook

$adminSheet = $wb.Worksheets.Item("Admin")
$adminSheet.Range("D2").Value = "1495934A"

$jiraSheet = $wb.Worksheets.Item("Jira")
$jiraSheet.Range("A3").Value = "Recruitment_RejectionList1"
